$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("capex_pv_calc")

# Remove the "MicroInverter Replacement" case-study rows (20-25), which
# contained a typo causing the B column to be 0 and the dependent formulas
# to resolve to #DIV/0!. Deleting the whole rows shifts everything below
# up by six rows and keeps downstream references (named formulas, etc.)
# consistent.
$ws.Range("A20:A25").EntireRow.Delete() | Out-Null

# The embedded chart is anchored by absolute row/column offsets that the
# host doesn't recompute automatically on a row delete, so nudge it up by
# the six rows (6 * 15pt row height) that just disappeared above it.
$co = $ws.ChartObjects(1)
$co.Top = $co.Top - 90
